$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C) for rows 2-28 from serial 45437 to 45439
for ($row = 2; $row -le 28; $row++) {
    $ws.Cells.Item($row, 3).Value = 45439
}
